# June newsletter - queue time
# Applies 4 changes:
#  1. Merge the "added two new visuals...caller." run with the following
#     single-space run into one run (same visible text, different run split).
#  2. Add the "Performance: On page 'Incoming Web Service performance', ..."
#     sentence (with a spell-check-split "analyze") to the empty paragraph
#     under "Performance report", and add a new "<insert screenshot>"
#     paragraph after it.
#  3. Move the <w:lastRenderedPageBreak/> from the second "Tips and Tricks"
#     heading run to the second "New signal/updates" heading run.

$d = $word.ActiveDocument

function Set-RangeOpenXml($targetRange, $innerXml) {
    $wrapperHead = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:xml="http://www.w3.org/XML/1998/namespace"><w:body>'
    $wrapperTail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $wrapped = $wrapperHead + $innerXml + $wrapperTail
    $targetRange.InsertXML($wrapped)
}

# ---------------------------------------------------------------------
# 1) Merge the two runs in the "IncomingWebserviceErrors" paragraph so
#    the trailing " " becomes part of the "...caller." run.
# ---------------------------------------------------------------------
$findRange1 = $d.Content
$found1 = $findRange1.Find.Execute("added two new visuals", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para1 = $findRange1.Paragraphs(1)
$para1Range = $para1.Range
$hunk1Xml = '<w:p w14:paraId="4C310745" w14:textId="64A9CDA2" w:rsidR="000940DF" w:rsidRPr="005E73D9" w:rsidRDefault="00E21918" w:rsidP="000940DF"><w:pPr><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:lastRenderedPageBreak/><w:t>On the page “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00E21918"><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t>IncomingWebserviceErrors</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t xml:space="preserve">”, </w:t></w:r><w:r w:rsidRPr="00E21918"><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t xml:space="preserve">added two new visuals. First one shows the new OData error dimensions ''Failure Reason'' and ''Diagnostic Message''. The second one shows the User Agent set by the caller. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t xml:space="preserve"> Use this to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t>a</w:t></w:r><w:r w:rsidRPr="00E21918"><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t>nalyze</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00E21918"><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t xml:space="preserve"> failed calls by user agents (who called). See failure reason and diagnostics messages (for API/OData requests) to help find the root cause of the errors.</w:t></w:r></w:p>'
Set-RangeOpenXml $para1Range $hunk1Xml

# ---------------------------------------------------------------------
# 2) Fill in the empty paragraph right after "Performance report" with
#    the queue-time sentence, then add a new "<insert screenshot>"
#    paragraph after it.
# ---------------------------------------------------------------------
$findRange2 = $d.Content
$found2 = $findRange2.Find.Execute("Performance report", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$perfPara = $findRange2.Paragraphs(1)
$targetPara = $perfPara.Next()
$targetParaRange = $targetPara.Range

$apos = [char]39
$perfSentence1 = '<w:p w14:paraId="0AB1D4EB" w14:textId="0239049D" w:rsidR="000940DF" w:rsidRPr="005E73D9" w:rsidRDefault="000940DF" w:rsidP="000940DF"><w:pPr><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t>Performance: On page '
$perfSentence2 = 'Incoming Web Service performance'
$perfSentence3 = ', added calculations on queue time (introduced in telemetry in v22)</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t>. Use this to a</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t>nalyze</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t xml:space="preserve"> time spent in web service queues</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$screenshotPara = '<w:p><w:pPr><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t>&lt;insert screenshot&gt;</w:t></w:r></w:p>'

$perfParaXml = $perfSentence1 + $apos + $perfSentence2 + $apos + $perfSentence3
$hunk2Xml = $perfParaXml + $screenshotPara

Set-RangeOpenXml $targetParaRange $hunk2Xml

# ---------------------------------------------------------------------
# 3) Move <w:lastRenderedPageBreak/> from the second "Tips and Tricks"
#    Heading-1 run to the second "New signal/updates" Heading-1 run.
#    (MatchCase=true + Style check so only the Heading-1 occurrence,
#    not the lower-case in-sentence mentions, is touched.)
# ---------------------------------------------------------------------
function Find-HeadingParagraph($searchText) {
    $rng = $d.Content
    while ($true) {
        $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $ok) {
            return $null
        }
        $candidate = $rng.Paragraphs(1)
        $styleName = $candidate.Style.NameLocal
        if ($styleName -eq "Heading 1") {
            return $candidate
        }
        $rng.Collapse(0)
    }
}

$newSignalPara = Find-HeadingParagraph "New signal/updates"
$newSignalParaRange = $newSignalPara.Range
$newSignalXml = '<w:p w14:paraId="13B041B9" w14:textId="3E418686" w:rsidR="0038500D" w:rsidRDefault="0038500D" w:rsidP="00D45DD4"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr></w:pPr><w:r w:rsidRPr="0038500D"><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:lastRenderedPageBreak/><w:t>New signal/updates</w:t></w:r></w:p>'
Set-RangeOpenXml $newSignalParaRange $newSignalXml

$tipsPara = Find-HeadingParagraph "Tips and Tricks"
$tipsParaRange = $tipsPara.Range
$tipsXml = '<w:p w14:paraId="3C722BA7" w14:textId="458A2910" w:rsidR="004D4487" w:rsidRPr="00615DEF" w:rsidRDefault="004D4487" w:rsidP="006B28FD"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr></w:pPr><w:r w:rsidRPr="00615DEF"><w:rPr><w:lang w:val="en-GB" w:eastAsia="en-DK"/></w:rPr><w:t>Tips and Tricks</w:t></w:r></w:p>'
Set-RangeOpenXml $tipsParaRange $tipsXml

Write-Output "Done: all 4 hunks applied."
